$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.820.04"
$ws.Range("E2").Value = "  +4.55%  "

$ws.Range("D3").Value = "2.630.18"
$ws.Range("E3").Value = "  +4.99%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'605.59"
$ws.Range("E5").Value = "  +1.64%  "

$ws.Range("D6").Value = "'178.83"
$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  +1.74%  "

$ws.Range("D9").Value = "2.633.09"
$ws.Range("E9").Value = "  +5.17%  "

$ws.Range("E10").Value = "  +8.20%  "

$ws.Range("D11").Value = "'0.166"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("E12").Value = "  +3.01%  "

$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").Value = "3.138.20"
$ws.Range("E14").Value = "  +5.98%  "

$ws.Range("E15").Value = "  +6.60%  "

$ws.Range("D16").Value = "72.678.74"
$ws.Range("E16").Value = "  +4.76%  "

$ws.Range("D17").Value = "'26.67"
$ws.Range("E17").Value = "  +3.08%  "

$ws.Range("D18").Value = "2.634.45"
$ws.Range("E18").Value = "  +4.61%  "

$ws.Range("D19").Value = "'384.91"
$ws.Range("E19").Value = "  +5.65%  "

$ws.Range("E20").Value = "  +6.71%  "

$ws.Range("E21").Value = "  +4.60%  "

$ws.Range("E22").Value = "  +3.25%  "

$ws.Range("E23").Value = "  +19.74%  "

$ws.Range("D24").Value = "'73.08"
$ws.Range("E24").Value = "  +3.64%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  +3.84%  "

$ws.Range("D27").Value = "'9.86"
$ws.Range("E27").Value = "  +8.97%  "

$ws.Range("D28").Value = "2.767.94"
$ws.Range("E28").Value = "  +4.76%  "

$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("D30").Value = "0.0₃0959"
$ws.Range("E30").Value = "  +7.28%  "

$ws.Range("D31").Value = "'533.26"
$ws.Range("E31").Value = "  +4.22%  "

$ws.Range("E32").Value = "  +3.69%  "

$ws.Range("E33").Value = "  +7.68%  "

$ws.Range("E34").Value = "  +3.15%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("D36").Value = "'163.44"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("D37").Value = "'19.32"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").Value = "'19.11"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.40"
$ws.Range("E39").Value = "  +7.02%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.111"
$ws.Range("E40").Value = "  -6.57%  "

$ws.Range("E41").Value = "  +5.43%  "

$ws.Range("D42").Value = "'5.07"
$ws.Range("E42").Value = "  +5.48%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  +12.31%  "

$ws.Range("E45").Value = "  +4.06%  "

$ws.Range("D46").Value = "'39.64"
$ws.Range("E46").Value = "  +2.03%  "

$ws.Range("D47").Value = "'151.06"
$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").Value = "'3.69"
$ws.Range("E48").Value = "  +3.05%  "

$ws.Range("E49").Value = "  +5.66%  "

$ws.Range("E50").Value = "  +8.20%  "

$ws.Range("E51").Value = "  +4.87%  "

Write-Host "edits applied"
